{"js": "// Replace the date heading and the division-problem cell values with the\n// new values from the target revision. Every \"find\" string in this map is\n// unique within the document, so a simple body-wide search & replace for\n// each pair is safe and unambiguous.\nconst replacements = [\n  [\"2025-06-17 Tuesday\", \"2025-06-18 Wednesday\"],\n  [\"43\u00f72=\", \"62\u00f79=\"],\n  [\"89\u00f76=\", \"25\u00f79=\"],\n  [\"12\u00f76=\", \"60\u00f76=\"],\n  [\"41\u00f72=\", \"79\u00f78=\"],\n  [\"57\u00f77=\", \"35\u00f77=\"],\n  [\"99\u00f78=\", \"67\u00f73=\"],\n  [\"36\u00f72=\", \"88\u00f75=\"],\n  [\"75\u00f72=\", \"18\u00f72=\"],\n  [\"97\u00f76=\", \"41\u00f79=\"],\n  [\"56\u00f74=\", \"23\u00f79=\"],\n  [\"77\u00f79=\", \"54\u00f78=\"],\n  [\"63\u00f76=\", \"71\u00f75=\"],\n  [\"34\u00f79=\", \"11\u00f75=\"],\n  [\"99\u00f79=\", \"96\u00f76=\"],\n  [\"19\u00f79=\", \"74\u00f78=\"],\n  [\"83\u00f79=\", \"66\u00f76=\"],\n  [\"15\u00f79=\", \"36\u00f77=\"],\n  [\"70\u00f79=\", \"90\u00f76=\"],\n  [\"68\u00f79=\", \"40\u00f72=\"],\n  [\"10\u00f77=\", \"12\u00f72=\"],\n  [\"29\u00f77=\", \"83\u00f78=\"],\n  [\"10\u00f72=\", \"76\u00f79=\"],\n  [\"28\u00f72=\", \"29\u00f75=\"],\n  [\"37\u00f76=\", \"87\u00f78=\"],\n  [\"24\u00f78=\", \"30\u00f76=\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and the division-problem cell values with the\n# new values from the target revision. Every \"find\" string below is unique\n# within the document, so Find/Replace (wrap = whole story, replace = all)\n# is safe and unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{find=\"2025-06-17 Tuesday\"; replace=\"2025-06-18 Wednesday\"},\n    @{find=\"43\u00f72=\"; replace=\"62\u00f79=\"},\n    @{find=\"89\u00f76=\"; replace=\"25\u00f79=\"},\n    @{find=\"12\u00f76=\"; replace=\"60\u00f76=\"},\n    @{find=\"41\u00f72=\"; replace=\"79\u00f78=\"},\n    @{find=\"57\u00f77=\"; replace=\"35\u00f77=\"},\n    @{find=\"99\u00f78=\"; replace=\"67\u00f73=\"},\n    @{find=\"36\u00f72=\"; replace=\"88\u00f75=\"},\n    @{find=\"75\u00f72=\"; replace=\"18\u00f72=\"},\n    @{find=\"97\u00f76=\"; replace=\"41\u00f79=\"},\n    @{find=\"56\u00f74=\"; replace=\"23\u00f79=\"},\n    @{find=\"77\u00f79=\"; replace=\"54\u00f78=\"},\n    @{find=\"63\u00f76=\"; replace=\"71\u00f75=\"},\n    @{find=\"34\u00f79=\"; replace=\"11\u00f75=\"},\n    @{find=\"99\u00f79=\"; replace=\"96\u00f76=\"},\n    @{find=\"19\u00f79=\"; replace=\"74\u00f78=\"},\n    @{find=\"83\u00f79=\"; replace=\"66\u00f76=\"},\n    @{find=\"15\u00f79=\"; replace=\"36\u00f77=\"},\n    @{find=\"70\u00f79=\"; replace=\"90\u00f76=\"},\n    @{find=\"68\u00f79=\"; replace=\"40\u00f72=\"},\n    @{find=\"10\u00f77=\"; replace=\"12\u00f72=\"},\n    @{find=\"29\u00f77=\"; replace=\"83\u00f78=\"},\n    @{find=\"10\u00f72=\"; replace=\"76\u00f79=\"},\n    @{find=\"28\u00f72=\"; replace=\"29\u00f75=\"},\n    @{find=\"37\u00f76=\"; replace=\"87\u00f78=\"},\n    @{find=\"24\u00f78=\"; replace=\"30\u00f76=\"}\n)\n\nforeach ($pair in $pairs) {\n    $range = $d.Content\n    $range.Find.Execute($pair.find, $true, $false, $false, $false, $false, $true, 1, $false, $pair.replace, 2) | Out-Null\n}\n"}
